$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18 (Leve Item ID 5471)
$ws.Cells.Item(18, 8).Value = 475
$ws.Cells.Item(18, 9).Value = 475
$ws.Cells.Item(18, 11).Value = 475
$ws.Cells.Item(18, 13).Value = -191

# Row 86 (Leve Item ID 12603)
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).ClearContents()
$ws.Cells.Item(86, 14).ClearContents()

# Row 89 (Leve Item ID 12603)
$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 13).ClearContents()
$ws.Cells.Item(89, 14).ClearContents()

# Row 95 (Leve Item ID 18200)
$ws.Cells.Item(95, 8).Value = 58600
$ws.Cells.Item(95, 10).Value = 58600
$ws.Cells.Item(95, 12).Value = 58600
$ws.Cells.Item(95, 14).Value = -64092

# Row 98 (Leve Item ID 36237)
$ws.Cells.Item(98, 8).Value = 1500
$ws.Cells.Item(98, 9).Value = 1170
$ws.Cells.Item(98, 10).Value = 2490
$ws.Cells.Item(98, 11).Value = 1170
$ws.Cells.Item(98, 12).Value = 2490
$ws.Cells.Item(98, 13).Value = 328
$ws.Cells.Item(98, 14).Value = -5486

# Row 111 (Leve Item ID 27768)
$ws.Cells.Item(111, 8).Value = 1000
$ws.Cells.Item(111, 9).Value = 1000
$ws.Cells.Item(111, 11).Value = 3000
$ws.Cells.Item(111, 13).Value = 67

# Row 122 (Leve Item ID 36237)
$ws.Cells.Item(122, 8).Value = 1500
$ws.Cells.Item(122, 9).Value = 1170
$ws.Cells.Item(122, 10).Value = 2490
$ws.Cells.Item(122, 11).Value = 3510
$ws.Cells.Item(122, 12).Value = 7470
$ws.Cells.Item(122, 13).Value = -1060
$ws.Cells.Item(122, 14).Value = -12370

$ws = $wb.Worksheets.Item("CRP")
# Row 17 (Leve Item ID 1823)
$ws.Cells.Item(17, 8).Value = 13945.5
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 13).ClearContents()

# Row 28 (Leve Item ID 18348)
$ws.Cells.Item(28, 8).Value = 15874.5
$ws.Cells.Item(28, 10).Value = 15874.5
$ws.Cells.Item(28, 12).Value = 15874.5
$ws.Cells.Item(28, 14).Value = -16364.5

# Row 92 (Leve Item ID 18041)
$ws.Cells.Item(92, 8).Value = 26766.834
$ws.Cells.Item(92, 10).Value = 28120.2
$ws.Cells.Item(92, 12).Value = 28120.2
$ws.Cells.Item(92, 14).Value = -33112.2

# Row 93 (Leve Item ID 19516)
$ws.Cells.Item(93, 8).Value = 12662.8
$ws.Cells.Item(93, 9).Value = 12662.8
$ws.Cells.Item(93, 11).Value = 12662.8
$ws.Cells.Item(93, 13).Value = -10790.8

# Row 95 (Leve Item ID 18192)
$ws.Cells.Item(95, 8).Value = 13806.5
$ws.Cells.Item(95, 10).Value = 13806.5
$ws.Cells.Item(95, 12).Value = 13806.5
$ws.Cells.Item(95, 14).Value = -19298.5

# Row 96 (Leve Item ID 18193)
$ws.Cells.Item(96, 8).Value = 14325.4
$ws.Cells.Item(96, 10).Value = 14325.4
$ws.Cells.Item(96, 12).Value = 14325.4
$ws.Cells.Item(96, 14).Value = -19817.4

# Row 99 (Leve Item ID 36198)
$ws.Cells.Item(99, 8).Value = 5154.154
$ws.Cells.Item(99, 10).Value = 5488.25
$ws.Cells.Item(99, 12).Value = 5488.25
$ws.Cells.Item(99, 14).Value = -8484.25

# Row 102 (Leve Item ID 19738)
$ws.Cells.Item(102, 8).Value = 40241
$ws.Cells.Item(102, 10).Value = 40241
$ws.Cells.Item(102, 12).Value = 40241
$ws.Cells.Item(102, 14).Value = -45109

# Row 103 (Leve Item ID 19558)
$ws.Cells.Item(103, 8).Value = 3762
$ws.Cells.Item(103, 9).Value = 3762
$ws.Cells.Item(103, 11).Value = 3762
$ws.Cells.Item(103, 13).Value = -2590

# Row 104 (Leve Item ID 19749)
$ws.Cells.Item(104, 8).Value = 40285
$ws.Cells.Item(104, 10).Value = 40285
$ws.Cells.Item(104, 12).Value = 40285
$ws.Cells.Item(104, 14).Value = -45527

# Row 106 (Leve Item ID 18661)
$ws.Cells.Item(106, 8).Value = 17671
$ws.Cells.Item(106, 10).Value = 17671
$ws.Cells.Item(106, 12).Value = 17671
$ws.Cells.Item(106, 14).Value = -20195

# Row 108 (Leve Item ID 27087)
$ws.Cells.Item(108, 8).Value = 79990
$ws.Cells.Item(108, 9).Value = 0
$ws.Cells.Item(108, 11).Value = 0
$ws.Cells.Item(108, 13).ClearContents()

# Row 109 (Leve Item ID 27203)
$ws.Cells.Item(109, 8).Value = 40220
$ws.Cells.Item(109, 10).Value = 40220
$ws.Cells.Item(109, 12).Value = 40220
$ws.Cells.Item(109, 14).Value = -42300

# Row 112 (Leve Item ID 25796)
$ws.Cells.Item(112, 8).Value = 0
$ws.Cells.Item(112, 10).Value = 0
$ws.Cells.Item(112, 12).Value = 0
$ws.Cells.Item(112, 14).ClearContents()

# Row 114 (Leve Item ID 27112)
$ws.Cells.Item(114, 8).Value = 50000
$ws.Cells.Item(114, 10).Value = 50000
$ws.Cells.Item(114, 12).Value = 50000
$ws.Cells.Item(114, 14).Value = -58678

# Row 117 (Leve Item ID 27135)
$ws.Cells.Item(117, 8).Value = 75000
$ws.Cells.Item(117, 9).Value = 75000
$ws.Cells.Item(117, 11).Value = 75000
$ws.Cells.Item(117, 13).Value = -70411

# Row 119 (Leve Item ID 26276)
$ws.Cells.Item(119, 8).Value = 50000
$ws.Cells.Item(119, 10).Value = 50000
$ws.Cells.Item(119, 12).Value = 50000
$ws.Cells.Item(119, 14).Value = -59676

# Row 125 (Leve Item ID 34297)
$ws.Cells.Item(125, 8).Value = 30000
$ws.Cells.Item(125, 10).Value = 30000
$ws.Cells.Item(125, 12).Value = 30000
$ws.Cells.Item(125, 14).Value = -34920

# Row 126 (Leve Item ID 36198)
$ws.Cells.Item(126, 8).Value = 5154.154
$ws.Cells.Item(126, 10).Value = 5488.25
$ws.Cells.Item(126, 12).Value = 16464.75
$ws.Cells.Item(126, 14).Value = -21404.75

$ws = $wb.Worksheets.Item("CUL")
# Row 37 (Leve Item ID 9516)
$ws.Cells.Item(37, 8).Value = 99648.336
$ws.Cells.Item(37, 10).Value = 99648.336
$ws.Cells.Item(37, 12).Value = 298945.008
$ws.Cells.Item(37, 14).Value = -299169.008

# Row 46 (Leve Item ID 4701)
$ws.Cells.Item(46, 8).Value = 1847.5
$ws.Cells.Item(46, 9).Value = 1075
$ws.Cells.Item(46, 10).Value = 2620
$ws.Cells.Item(46, 11).Value = 3225
$ws.Cells.Item(46, 12).Value = 7860
$ws.Cells.Item(46, 13).Value = -3134
$ws.Cells.Item(46, 14).Value = -8042

$ws = $wb.Worksheets.Item("GSM")
# Row 44 (Leve Item ID 4143)
$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(44, 10).Value = 0
$ws.Cells.Item(44, 12).Value = 0
$ws.Cells.Item(44, 14).ClearContents()

# Row 102 (Leve Item ID 36169)
$ws.Cells.Item(102, 8).Value = 2372.5
$ws.Cells.Item(102, 9).Value = 1830.3334
$ws.Cells.Item(102, 11).Value = 1830.3334
$ws.Cells.Item(102, 13).Value = -208.3334

# Row 126 (Leve Item ID 36184)
$ws.Cells.Item(126, 8).Value = 15637.5
$ws.Cells.Item(126, 9).Value = 14952.75
$ws.Cells.Item(126, 10).Value = 17007
$ws.Cells.Item(126, 11).Value = 44858.25
$ws.Cells.Item(126, 12).Value = 51021
$ws.Cells.Item(126, 13).Value = -42388.25
$ws.Cells.Item(126, 14).Value = -55961

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Cells.Item(7, 8).Value = 27701.908
$ws.Cells.Item(7, 9).Value = 27747.1
$ws.Cells.Item(7, 11).Value = 27747.1
$ws.Cells.Item(7, 13).Value = -27635.1

# Row 22 (Leve Item ID 5277)
$ws.Cells.Item(22, 8).Value = 1015.4375
$ws.Cells.Item(22, 9).Value = 816.2222
$ws.Cells.Item(22, 10).Value = 1271.5714
$ws.Cells.Item(22, 11).Value = 816.2222
$ws.Cells.Item(22, 12).Value = 1271.5714
$ws.Cells.Item(22, 13).Value = -521.2222
$ws.Cells.Item(22, 14).Value = -1861.5714

# Row 27 (Leve Item ID 5277)
$ws.Cells.Item(27, 8).Value = 1015.4375
$ws.Cells.Item(27, 9).Value = 816.2222
$ws.Cells.Item(27, 10).Value = 1271.5714
$ws.Cells.Item(27, 11).Value = 816.2222
$ws.Cells.Item(27, 12).Value = 1271.5714
$ws.Cells.Item(27, 13).Value = -709.2222
$ws.Cells.Item(27, 14).Value = -1485.5714

# Row 40 (Leve Item ID 36248)
$ws.Cells.Item(40, 8).Value = 7710
$ws.Cells.Item(40, 9).Value = 7710
$ws.Cells.Item(40, 11).Value = 7710
$ws.Cells.Item(40, 13).Value = -7574

# Row 46 (Leve Item ID 5282)
$ws.Cells.Item(46, 8).Value = 3875.25
$ws.Cells.Item(46, 10).Value = 4750.5
$ws.Cells.Item(46, 12).Value = 4750.5
$ws.Cells.Item(46, 14).Value = -5126.5

# Row 122 (Leve Item ID 36247)
$ws.Cells.Item(122, 8).Value = 3504
$ws.Cells.Item(122, 9).Value = 3504
$ws.Cells.Item(122, 11).Value = 10512
$ws.Cells.Item(122, 13).Value = -8062

# Row 126 (Leve Item ID 36249)
$ws.Cells.Item(126, 8).Value = 27701.908
$ws.Cells.Item(126, 9).Value = 27747.1
$ws.Cells.Item(126, 11).Value = 83241.29999999999
$ws.Cells.Item(126, 13).Value = -80771.29999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (Leve Item ID 44029)
$ws.Cells.Item(132, 8).Value = 888.5
$ws.Cells.Item(132, 9).Value = 888.5
$ws.Cells.Item(132, 11).Value = 2665.5
$ws.Cells.Item(132, 13).Value = -135.5
